$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.871.73"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "2.083.98"
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.29"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.26"
$ws.Range("E7").Value = "  +3.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.396"
$ws.Range("E9").Value = "  +2.10%  "

$ws.Range("E10").Value = "  +1.42%  "

$ws.Range("E11").Value = "  +1.25%  "

$ws.Range("E12").Value = "  +2.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.18"
$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.776"
$ws.Range("E14").Value = "  +1.53%  "

$ws.Range("E15").Value = "  +2.21%  "

$ws.Range("D16").Value = "2.061.13"
$ws.Range("E16").Value = "  -1.88%  "

$ws.Range("D17").Value = "37.772.51"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.61"
$ws.Range("E19").Value = "  +0.95%  "

$ws.Range("E20").Value = "  +3.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.14"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("E23").Value = "  -0.70%  "

$ws.Range("E24").Value = "  +1.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.99"
$ws.Range("E25").Value = "  +2.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.16"
$ws.Range("E26").Value = "  +2.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.137"
$ws.Range("E27").Value = "  -2.70%  "

$ws.Range("E28").Value = "  -1.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.50"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("E30").Value = "  +1.86%  "

$ws.Range("E31").Value = "  +2.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.74"

$ws.Range("E33").Value = "  +1.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.51"
$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.45"
$ws.Range("E35").Value = "  +1.84%  "

$ws.Range("E36").Value = "  -0.63%  "

$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.42"
$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0982"
$ws.Range("E39").Value = "  -1.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.18"
$ws.Range("E40").Value = "  +2.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0218"
$ws.Range("E41").Value = "  +1.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.90"
$ws.Range("E42").Value = "  -1.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.80"
$ws.Range("E43").Value = "  +7.40%  "

$ws.Range("D44").Value = "1.448.20"
$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.16"
$ws.Range("E45").Value = "  -0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.19"
$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("E47").Value = "  +0.68%  "

$ws.Range("E48").Value = "  +0.65%  "

$ws.Range("E49").Value = "  -0.43%  "

$ws.Range("D50").Value = "2.274.81"
$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.95"
$ws.Range("E51").Value = "  +1.15%  "
